$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.807.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.287.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.37%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "265.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.21"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.77%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.91%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.54%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.95"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.11%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.625.95"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.61%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.824"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.251.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.835.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +13.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.25"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.05%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.86%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.44"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.07%  "

# Row 28
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.36"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.31%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.86%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.08"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.91%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.71%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.50"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.50%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.80%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.37"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0343"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.78%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +13.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.240"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +19.87%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.29"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.37%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.97%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.42"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.78%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.27"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +13.44%  "

# Row 45
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.55"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.70%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.36%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.96"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.65%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.27%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.508.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.18%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.423"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.58%  "
